# Generate Report for Handback
#
# The 425d36a0-... file has now been handed back (in addition to the
# a9377a48-... file that was already handed back).  This script updates
# the localization-status report:
#   - Overview sheet: both files now show "Handed back: in sync with en-US"
#   - zh-cn / de-de sheets: both files now show status "Handed back" and a
#     new "Latest Handback DateTime" is recorded for the 425d file.
# Because the rows are keyed by file, and the 425d row now sorts before
# the a9377 row (matching the authoritative export order), the two data
# rows on every sheet swap their content.

$wb = $excel.ActiveWorkbook

$A9377 = "a9377a48-67ad-4e43-b924-f3c40a14ff5b"
$D425  = "425d36a0-c691-45e2-ab11-0a37ec306df5"

$statusHandedBack = "Handed back: in sync with en-US"

function Set-CellValues {
    param($ws, $values)
    foreach ($ref in $values.Keys) {
        $ws.Range($ref).Value = $values[$ref]
    }
}

function Set-HyperlinkDisplay {
    param($ws, $displays)
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($displays.ContainsKey($addr)) {
            $hl.TextToDisplay = $displays[$addr]
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellValues $wsOverview @{
    "A2" = "$D425.md"
    "B2" = $statusHandedBack
    "C2" = $statusHandedBack
    "D2" = "2016-37-17 18:37:55"

    "A3" = "$A9377.md"
    "B3" = $statusHandedBack
    "C3" = $statusHandedBack
    "D3" = "2016-36-17 18:36:44"
}

Set-HyperlinkDisplay $wsOverview @{
    '$A$2' = "$D425.md"
    '$A$3' = "$A9377.md"
}

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status |
#   Latest Handoff File | Latest Handoff Datetime | Latest Target File |
#   Latest Handback File | Latest Handback DateTime | Handoff Reason
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellValues $wsZhCn @{
    "A2" = "$D425.md"
    "B2" = ".md"
    "C2" = $statusHandedBack
    "D2" = "$D425.30a09fc566298713cba6fd8bbaf67821415842ff.zh-cn.xlf"
    "E2" = "2016-03-17 18:37:52"
    "F2" = "$D425.md"
    "G2" = "$D425.30a09fc566298713cba6fd8bbaf67821415842ff.zh-cn.xlf"
    "H2" = "2016-03-17 18:38:11"
    "I2" = "Include"

    "A3" = "$A9377.md"
    "B3" = ".md"
    "C3" = $statusHandedBack
    "D3" = "$A9377.42150e9487e50f23fbeb2f83d65c9a386e6f9760.zh-cn.xlf"
    "E3" = "2016-03-17 18:36:29"
    "F3" = "$A9377.md"
    "G3" = "$A9377.42150e9487e50f23fbeb2f83d65c9a386e6f9760.zh-cn.xlf"
    "H3" = "2016-03-17 18:37:23"
    "I3" = "Include"
}

Set-HyperlinkDisplay $wsZhCn @{
    '$A$2' = "$D425.md"
    '$B$2' = ".md"
    '$D$2' = "$D425.30a09fc566298713cba6fd8bbaf67821415842ff.zh-cn.xlf"
    '$F$2' = "$D425.md"
    '$G$2' = "$D425.30a09fc566298713cba6fd8bbaf67821415842ff.zh-cn.xlf"

    '$A$3' = "$A9377.md"
    '$B$3' = ".md"
    '$D$3' = "$A9377.42150e9487e50f23fbeb2f83d65c9a386e6f9760.zh-cn.xlf"
    '$F$3' = "$A9377.md"
    '$G$3' = "$A9377.42150e9487e50f23fbeb2f83d65c9a386e6f9760.zh-cn.xlf"
}

# ---------------------------------------------------------------------
# Sheet "de-de": same layout as zh-cn
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellValues $wsDeDe @{
    "A2" = "$D425.md"
    "B2" = ".md"
    "C2" = $statusHandedBack
    "D2" = "$D425.30a09fc566298713cba6fd8bbaf67821415842ff.de-de.xlf"
    "E2" = "2016-03-17 18:37:55"
    "F2" = "$D425.md"
    "G2" = "$D425.30a09fc566298713cba6fd8bbaf67821415842ff.de-de.xlf"
    "H2" = "2016-03-17 18:38:19"
    "I2" = "Include"

    "A3" = "$A9377.md"
    "B3" = ".md"
    "C3" = $statusHandedBack
    "D3" = "$A9377.42150e9487e50f23fbeb2f83d65c9a386e6f9760.de-de.xlf"
    "E3" = "2016-03-17 18:36:44"
    "F3" = "$A9377.md"
    "G3" = "$A9377.42150e9487e50f23fbeb2f83d65c9a386e6f9760.de-de.xlf"
    "H3" = "2016-03-17 18:37:29"
    "I3" = "Include"
}

Set-HyperlinkDisplay $wsDeDe @{
    '$A$2' = "$D425.md"
    '$B$2' = ".md"
    '$D$2' = "$D425.30a09fc566298713cba6fd8bbaf67821415842ff.de-de.xlf"
    '$F$2' = "$D425.md"
    '$G$2' = "$D425.30a09fc566298713cba6fd8bbaf67821415842ff.de-de.xlf"

    '$A$3' = "$A9377.md"
    '$B$3' = ".md"
    '$D$3' = "$A9377.42150e9487e50f23fbeb2f83d65c9a386e6f9760.de-de.xlf"
    '$F$3' = "$A9377.md"
    '$G$3' = "$A9377.42150e9487e50f23fbeb2f83d65c9a386e6f9760.de-de.xlf"
}
